$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 499.83334
$ws.Range("I18").Value = 559.4
$ws.Range("J18").Value = 202
$ws.Range("K18").Value = 559.4
$ws.Range("L18").Value = 202
$ws.Range("M18").Value = -275.4
$ws.Range("N18").Value = -770
$ws.Range("H19").Value = 6099.75
$ws.Range("I19").Value = 7999
$ws.Range("J19").Value = 5466.6665
$ws.Range("K19").Value = 7999
$ws.Range("L19").Value = 5466.6665
$ws.Range("M19").Value = -7824
$ws.Range("N19").Value = -5816.6665
$ws.Range("H51").Value = 85419560
$ws.Range("J51").Value = 3275
$ws.Range("L51").Value = 3275
$ws.Range("N51").Value = -4243
$ws.Range("H132").Value = 4893.2095
$ws.Range("I132").Value = 3262.3948
$ws.Range("K132").Value = 9787.1844
$ws.Range("M132").Value = -7257.1844
$ws.Range("H137").Value = 2506.5652
$ws.Range("J137").Value = 2445.25
$ws.Range("L137").Value = 7335.75
$ws.Range("N137").Value = -12435.75
$ws.Range("H138").Value = 276532.7
$ws.Range("I138").Value = 34652.266
$ws.Range("J138").Value = 1002173.9
$ws.Range("K138").Value = 103956.798
$ws.Range("L138").Value = 3006521.7
$ws.Range("M138").Value = -98816.79800000001
$ws.Range("N138").Value = -3016801.7
$ws.Range("H139").Value = 85000
$ws.Range("J139").Value = 85000
$ws.Range("L139").Value = 85000
$ws.Range("N139").Value = -95280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1320.9062
$ws.Range("I2").Value = 1290.6296
$ws.Range("K2").Value = 1290.6296
$ws.Range("M2").Value = -1177.6296
$ws.Range("H116").Value = 1320.9062
$ws.Range("I116").Value = 1290.6296
$ws.Range("K116").Value = 1290.6296
$ws.Range("M116").Value = 1003.3704
$ws.Range("H122").Value = 1792.025
$ws.Range("I122").Value = 1618.0541
$ws.Range("K122").Value = 4854.1623
$ws.Range("M122").Value = -2404.1623

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1320.9062
$ws.Range("I3").Value = 1290.6296
$ws.Range("K3").Value = 1290.6296
$ws.Range("M3").Value = -1176.6296
$ws.Range("H33").Value = 12000
$ws.Range("I33").Value = 12000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 12000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -11664
$ws.Range("N33").Value = $null
$ws.Range("H44").Value = 25000
$ws.Range("J44").Value = 25000
$ws.Range("L44").Value = 25000
$ws.Range("N44").Value = -25994
$ws.Range("H49").Value = 19333
$ws.Range("J49").Value = 19333
$ws.Range("L49").Value = 19333
$ws.Range("N49").Value = -19811
$ws.Range("H134").Value = 6175.93
$ws.Range("I134").Value = 2431.476
$ws.Range("K134").Value = 7294.428
$ws.Range("M134").Value = -4759.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 299.5238
$ws.Range("I7").Value = 257.83334
$ws.Range("J7").Value = 549.6667
$ws.Range("K7").Value = 257.83334
$ws.Range("L7").Value = 549.6667
$ws.Range("M7").Value = -144.83334
$ws.Range("N7").Value = -775.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 146666
$ws.Range("I128").Value = 146666
$ws.Range("K128").Value = 439998
$ws.Range("M128").Value = -435018

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 281.9
$ws.Range("I2").Value = 216
$ws.Range("J2").Value = 347.8
$ws.Range("K2").Value = 216
$ws.Range("L2").Value = 347.8
$ws.Range("M2").Value = -103
$ws.Range("N2").Value = -573.8
$ws.Range("H6").Value = 1208
$ws.Range("J6").Value = 1208
$ws.Range("L6").Value = 1208
$ws.Range("N6").Value = -1434
$ws.Range("H16").Value = 1208
$ws.Range("J16").Value = 1208
$ws.Range("L16").Value = 1208
$ws.Range("N16").Value = -1708
$ws.Range("H17").Value = 3752.6667
$ws.Range("J17").Value = 625
$ws.Range("L17").Value = 625
$ws.Range("N17").Value = -961
$ws.Range("H18").Value = 31111
$ws.Range("J18").Value = 31111
$ws.Range("L18").Value = 31111
$ws.Range("N18").Value = -31697
$ws.Range("H19").Value = 23950
$ws.Range("I19").Value = 12000
$ws.Range("K19").Value = 12000
$ws.Range("M19").Value = -11712
$ws.Range("H34").Value = 166831.4
$ws.Range("J34").Value = 180000
$ws.Range("L34").Value = 180000
$ws.Range("N34").Value = -180536
$ws.Range("H39").Value = 80073.86
$ws.Range("J39").Value = 80073.86
$ws.Range("L39").Value = 80073.86
$ws.Range("N39").Value = -81137.86
$ws.Range("H76").Value = 166831.4
$ws.Range("J76").Value = 180000
$ws.Range("L76").Value = 180000
$ws.Range("N76").Value = -180630
$ws.Range("H79").Value = 166831.4
$ws.Range("J79").Value = 180000
$ws.Range("L79").Value = 180000
$ws.Range("N79").Value = -182184
$ws.Range("H122").Value = 1924.5769
$ws.Range("I122").Value = 1697.3914
$ws.Range("K122").Value = 5092.174199999999
$ws.Range("M122").Value = -2642.174199999999
$ws.Range("H126").Value = 2842.353
$ws.Range("I126").Value = 2831
$ws.Range("J126").Value = 2927.5
$ws.Range("K126").Value = 8493
$ws.Range("L126").Value = 8782.5
$ws.Range("M126").Value = -6023
$ws.Range("N126").Value = -13722.5
$ws.Range("H132").Value = 14497458
$ws.Range("I132").Value = 17548566
$ws.Range("J132").Value = 4693.75
$ws.Range("K132").Value = 52645698
$ws.Range("L132").Value = 14081.25
$ws.Range("M132").Value = -52643168
$ws.Range("N132").Value = -19141.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 10445.8
$ws.Range("I46").Value = 10380
$ws.Range("J46").Value = 10455.923
$ws.Range("K46").Value = 10380
$ws.Range("L46").Value = 10455.923
$ws.Range("M46").Value = -10192
$ws.Range("N46").Value = -10831.923
$ws.Range("H97").Value = 60421.75
$ws.Range("J97").Value = 60421.75
$ws.Range("L97").Value = 60421.75
$ws.Range("N97").Value = -62403.75
$ws.Range("H122").Value = 4469.6294
$ws.Range("I122").Value = 3750.8096
$ws.Range("K122").Value = 11252.4288
$ws.Range("M122").Value = -8802.4288
$ws.Range("H136").Value = 3752.6943
$ws.Range("I136").Value = 3570.742
$ws.Range("K136").Value = 10712.226
$ws.Range("M136").Value = -8162.226000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 3600.6667
$ws.Range("I17").Value = 3600.6667
$ws.Range("K17").Value = 3600.6667
$ws.Range("M17").Value = -3428.6667
$ws.Range("H18").Value = 333346660
$ws.Range("J18").Value = 333346660
$ws.Range("L18").Value = 333346660
$ws.Range("N18").Value = -333347006
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = $null
$ws.Range("H113").Value = 2874384.5
$ws.Range("I113").Value = 5555890.5
$ws.Range("J113").Value = 1342.5
$ws.Range("K113").Value = 16667671.5
$ws.Range("L113").Value = 4027.5
$ws.Range("M113").Value = -16665501.5
$ws.Range("N113").Value = -8367.5
